$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 365.38
$ws.Range("D2").Value = 335.55

$ws.Range("C3").Value = 123.15
$ws.Range("D3").Value = 141.41

$ws.Range("C4").Value = 128.3
$ws.Range("D4").Value = 89.16

$ws.Range("C5").Value = 113.93
$ws.Range("D5").Value = 104.98

$ws.Range("C6").Value = 67.31
$ws.Range("D6").Value = 64

$ws.Range("C7").Value = 73.20999999999999
$ws.Range("D7").Value = 45.25

$ws.Range("C8").Value = 55.84
$ws.Range("D8").Value = 77.41

$ws.Range("C9").Value = 66.41
$ws.Range("D9").Value = 52.57

$ws.Range("C10").Value = 61.9
$ws.Range("D10").Value = 36.58

$ws.Range("C11").Value = 40.72
$ws.Range("D11").Value = 59.73

$wb.Save()
